$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("SmaCross") gets recomputed backtest stats (output re-run with a
# different strategy/parameter choice).
$ws.Range("B3").Value = 504856.049390626
$ws.Range("C3").Value = 1010532.421640625
$ws.Range("D3").Value = 399923.9929531248
$ws.Range("E3").Value = -49.5143950609374
$ws.Range("F3").Value = -4.25731054355469
$ws.Range("G3").Value = 156
$ws.Range("H3").Value = 34.61538461538461
$ws.Range("I3").Value = 6.721220520353199
$ws.Range("J3").Value = -4.439136829842005
$ws.Range("K3").Value = -0.06925583491882437
$ws.Range("L3").Value = 1.270833333333333
$ws.Range("M3").Value = 0.1958333333333333

# The duplicate SmaCross row and the extra EmaCross row are no longer part
# of the output selection, so remove rows 4 and 5 entirely (sheet shrinks
# from A1:M5 down to A1:M3).
$ws.Rows("4:5").Delete()
